$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45
$ws.Range("B45").Value = 6803331
$ws.Range("E45").Value = "Hanacka Slavia Kromeriz"
$ws.Range("F45").Value = "Opava"
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 1
$ws.Range("I45").Value = "A"
$ws.Range("J45").Value = 2.15
$ws.Range("K45").Value = 3.2
$ws.Range("L45").Value = 2.9
$ws.Range("M45").Value = 2.5
$ws.Range("N45").Value = 3.2
$ws.Range("O45").Value = 2.5
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 1.95
$ws.Range("R45").Value = 1.85
$ws.Range("S45").Value = 2.5
$ws.Range("T45").Value = 1.975
$ws.Range("U45").Value = 1.825
$ws.Range("V45").Value = -1
$ws.Range("W45").Value = -1
$ws.Range("X45").Value = 1.5
$ws.Range("Y45").Value = -1
$ws.Range("Z45").Value = 0.8500000000000001
$ws.Range("AA45").Value = -1
$ws.Range("AB45").Value = 0.825

# Row 46
$ws.Range("B46").Value = 6804171
$ws.Range("E46").Value = "SK Lisen"
$ws.Range("F46").Value = "Varnsdorf"
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 1
$ws.Range("I46").Value = "H"
$ws.Range("J46").Value = 1.909
$ws.Range("K46").Value = 3.3
$ws.Range("L46").Value = 3.5
$ws.Range("M46").Value = 1.909
$ws.Range("N46").Value = 3.3
$ws.Range("O46").Value = 3.5
$ws.Range("P46").Value = -0.5
$ws.Range("Q46").Value = 2
$ws.Range("R46").Value = 1.8
$ws.Range("S46").Value = 2.5
$ws.Range("T46").Value = 1.85
$ws.Range("U46").Value = 1.95
$ws.Range("V46").Value = 0.909
$ws.Range("W46").Value = -1
$ws.Range("X46").Value = -1
$ws.Range("Y46").Value = 1
$ws.Range("Z46").Value = -1
$ws.Range("AA46").Value = 0.8500000000000001
$ws.Range("AB46").Value = -1

# Row 47
$ws.Range("B47").Value = 6804172
$ws.Range("E47").Value = "Dukla Praha"
$ws.Range("F47").Value = "Vysocina Jihlava"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = "H"
$ws.Range("J47").Value = 1.8
$ws.Range("K47").Value = 3.4
$ws.Range("L47").Value = 3.8
$ws.Range("M47").Value = 1.5
$ws.Range("N47").Value = 4.2
$ws.Range("O47").Value = 5
$ws.Range("P47").Value = -1
$ws.Range("Q47").Value = 1.825
$ws.Range("R47").Value = 1.975
$ws.Range("S47").Value = 3
$ws.Range("T47").Value = 1.925
$ws.Range("U47").Value = 1.875
$ws.Range("V47").Value = 0.5
$ws.Range("W47").Value = -1
$ws.Range("X47").Value = -1
$ws.Range("Y47").Value = 0
$ws.Range("Z47").Value = 0
$ws.Range("AA47").Value = -1
$ws.Range("AB47").Value = 0.875

# Row 48
$ws.Range("B48").Value = 6804174
$ws.Range("E48").Value = "FC Silon Taborsko"
$ws.Range("F48").Value = "MFK Vyskov"
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = "D"
$ws.Range("J48").Value = 2.4
$ws.Range("K48").Value = 3.2
$ws.Range("L48").Value = 2.625
$ws.Range("M48").Value = 2.45
$ws.Range("N48").Value = 3.25
$ws.Range("O48").Value = 2.625
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 1.85
$ws.Range("R48").Value = 1.95
$ws.Range("S48").Value = 2.75
$ws.Range("T48").Value = 1.975
$ws.Range("U48").Value = 1.825
$ws.Range("V48").Value = -1
$ws.Range("W48").Value = 2.25
$ws.Range("X48").Value = -1
$ws.Range("Y48").Value = 0
$ws.Range("Z48").Value = 0
$ws.Range("AA48").Value = -1
$ws.Range("AB48").Value = 0.825

# Row 49
$ws.Range("B49").Value = 6804173
$ws.Range("E49").Value = "FC Vlasim"
$ws.Range("F49").Value = "FK Pribram"
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1
$ws.Range("I49").Value = "H"
$ws.Range("J49").Value = 2.25
$ws.Range("K49").Value = 3.25
$ws.Range("L49").Value = 2.75
$ws.Range("M49").Value = 1.833
$ws.Range("N49").Value = 3.8
$ws.Range("O49").Value = 3.3
$ws.Range("P49").Value = -0.5
$ws.Range("Q49").Value = 1.875
$ws.Range("R49").Value = 1.925
$ws.Range("S49").Value = 3
$ws.Range("T49").Value = 1.875
$ws.Range("U49").Value = 1.925
$ws.Range("V49").Value = 0.833
$ws.Range("W49").Value = -1
$ws.Range("X49").Value = -1
$ws.Range("Y49").Value = 0.875
$ws.Range("Z49").Value = -1
$ws.Range("AA49").Value = 0.875
$ws.Range("AB49").Value = -1

# Row 60
$ws.Range("B60").Value = 6803164
$ws.Range("E60").Value = "FC Brno"
$ws.Range("F60").Value = "Varnsdorf"
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = "H"
$ws.Range("J60").Value = 1.6
$ws.Range("K60").Value = 3.75
$ws.Range("L60").Value = 4.5
$ws.Range("M60").Value = 1.363
$ws.Range("N60").Value = 4.333
$ws.Range("O60").Value = 6.5
$ws.Range("P60").Value = -1.25
$ws.Range("Q60").Value = 1.8
$ws.Range("R60").Value = 2
$ws.Range("S60").Value = 3
$ws.Range("T60").Value = 1.975
$ws.Range("U60").Value = 1.825
$ws.Range("V60").Value = 0.363
$ws.Range("W60").Value = -1
$ws.Range("X60").Value = -1
$ws.Range("Y60").Value = -0.5
$ws.Range("Z60").Value = 0.5
$ws.Range("AA60").Value = -1
$ws.Range("AB60").Value = 0.825

# Row 61
$ws.Range("B61").Value = 6803246
$ws.Range("E61").Value = "Hanacka Slavia Kromeriz"
$ws.Range("F61").Value = "Viktoria Zizkov"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = "H"
$ws.Range("J61").Value = 2.5
$ws.Range("K61").Value = 3.2
$ws.Range("L61").Value = 2.5
$ws.Range("M61").Value = 2.7
$ws.Range("N61").Value = 3.3
$ws.Range("O61").Value = 2.3
$ws.Range("P61").Value = 0
$ws.Range("Q61").Value = 1.95
$ws.Range("R61").Value = 1.85
$ws.Range("S61").Value = 2.75
$ws.Range("T61").Value = 2
$ws.Range("U61").Value = 1.8
$ws.Range("V61").Value = 1.7
$ws.Range("W61").Value = -1
$ws.Range("X61").Value = -1
$ws.Range("Y61").Value = 0.95
$ws.Range("Z61").Value = -1
$ws.Range("AA61").Value = -1
$ws.Range("AB61").Value = 0.8

# Row 64
$ws.Range("B64").Value = 6804181
$ws.Range("E64").Value = "Sigma Olomouc B"
$ws.Range("F64").Value = "FK Pribram"
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 1
$ws.Range("I64").Value = "H"
$ws.Range("J64").Value = 2.4
$ws.Range("K64").Value = 3.25
$ws.Range("L64").Value = 2.55
$ws.Range("M64").Value = 2
$ws.Range("N64").Value = 3.3
$ws.Range("O64").Value = 3.2
$ws.Range("P64").Value = -0.25
$ws.Range("Q64").Value = 1.8
$ws.Range("R64").Value = 2
$ws.Range("S64").Value = 2.5
$ws.Range("T64").Value = 1.825
$ws.Range("U64").Value = 1.975
$ws.Range("V64").Value = 1
$ws.Range("W64").Value = -1
$ws.Range("X64").Value = -1
$ws.Range("Y64").Value = 0.8
$ws.Range("Z64").Value = -1
$ws.Range("AA64").Value = 0.825
$ws.Range("AB64").Value = -1

# Row 65
$ws.Range("B65").Value = 6804185
$ws.Range("E65").Value = "Sparta Prague B"
$ws.Range("F65").Value = "Opava"
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = "D"
$ws.Range("J65").Value = 2
$ws.Range("K65").Value = 3.4
$ws.Range("L65").Value = 3.1
$ws.Range("M65").Value = 2.625
$ws.Range("N65").Value = 3.5
$ws.Range("O65").Value = 2.25
$ws.Range("P65").Value = 0.25
$ws.Range("Q65").Value = 1.75
$ws.Range("R65").Value = 2.05
$ws.Range("S65").Value = 2.75
$ws.Range("T65").Value = 1.975
$ws.Range("U65").Value = 1.825
$ws.Range("V65").Value = -1
$ws.Range("W65").Value = 2.5
$ws.Range("X65").Value = -1
$ws.Range("Y65").Value = 0.375
$ws.Range("Z65").Value = -0.5
$ws.Range("AA65").Value = -1
$ws.Range("AB65").Value = 0.825

# Row 124
$ws.Range("B124").Value = 6804231
$ws.Range("E124").Value = "MFK Chrudim"
$ws.Range("F124").Value = "FC Vlasim"
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = "H"
$ws.Range("J124").Value = 2.35
$ws.Range("K124").Value = 3.25
$ws.Range("L124").Value = 2.6
$ws.Range("M124").Value = 2.45
$ws.Range("N124").Value = 3.3
$ws.Range("O124").Value = 2.45
$ws.Range("P124").Value = 0
$ws.Range("Q124").Value = 1.925
$ws.Range("R124").Value = 1.875
$ws.Range("S124").Value = 2.75
$ws.Range("T124").Value = 2
$ws.Range("U124").Value = 1.8
$ws.Range("V124").Value = 1.45
$ws.Range("W124").Value = -1
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = 0.925
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = 0.8

# Row 125
$ws.Range("B125").Value = 6804230
$ws.Range("E125").Value = "MFK Vyskov"
$ws.Range("F125").Value = "Dukla Praha"
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = "H"
$ws.Range("J125").Value = 2.05
$ws.Range("K125").Value = 3.25
$ws.Range("L125").Value = 3.1
$ws.Range("M125").Value = 1.95
$ws.Range("N125").Value = 3.2
$ws.Range("O125").Value = 3.4
$ws.Range("P125").Value = -0.5
$ws.Range("Q125").Value = 2
$ws.Range("R125").Value = 1.8
$ws.Range("S125").Value = 2.5
$ws.Range("T125").Value = 1.975
$ws.Range("U125").Value = 1.825
$ws.Range("V125").Value = 0.95
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 1
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = -1
$ws.Range("AB125").Value = 0.825

# Row 155
$ws.Range("B155").Value = 6804251
$ws.Range("E155").Value = "FK Pribram"
$ws.Range("F155").Value = "FC Sellier  Bellot Vlasim"
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 1
$ws.Range("I155").Value = "A"
$ws.Range("J155").Value = 2.1
$ws.Range("K155").Value = 3.25
$ws.Range("L155").Value = 3
$ws.Range("M155").Value = 1.909
$ws.Range("N155").Value = 3.5
$ws.Range("O155").Value = 3.5
$ws.Range("P155").Value = -0.5
$ws.Range("Q155").Value = 1.975
$ws.Range("R155").Value = 1.825
$ws.Range("S155").Value = 2.5
$ws.Range("T155").Value = 2
$ws.Range("U155").Value = 1.8
$ws.Range("V155").Value = -1
$ws.Range("W155").Value = -1
$ws.Range("X155").Value = 2.5
$ws.Range("Y155").Value = -1
$ws.Range("Z155").Value = 0.825
$ws.Range("AA155").Value = -1
$ws.Range("AB155").Value = 0.8

# Row 156
$ws.Range("B156").Value = 6804250
$ws.Range("E156").Value = "Vysocina Jihlava"
$ws.Range("F156").Value = "Dukla Praha"
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 2
$ws.Range("I156").Value = "A"
$ws.Range("J156").Value = 2.75
$ws.Range("K156").Value = 3.4
$ws.Range("L156").Value = 2.2
$ws.Range("M156").Value = 3.1
$ws.Range("N156").Value = 3.4
$ws.Range("O156").Value = 2.1
$ws.Range("P156").Value = 0.25
$ws.Range("Q156").Value = 1.95
$ws.Range("R156").Value = 1.85
$ws.Range("S156").Value = 2.5
$ws.Range("T156").Value = 1.8
$ws.Range("U156").Value = 2
$ws.Range("V156").Value = -1
$ws.Range("W156").Value = -1
$ws.Range("X156").Value = 1.1
$ws.Range("Y156").Value = -1
$ws.Range("Z156").Value = 0.8500000000000001
$ws.Range("AA156").Value = 0.8
$ws.Range("AB156").Value = -1

# Row 170
$ws.Range("B170").Value = 6804259
$ws.Range("E170").Value = "FK Pribram"
$ws.Range("F170").Value = "Sigma Olomouc B"
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 1
$ws.Range("I170").Value = "A"
$ws.Range("J170").Value = 1.952
$ws.Range("K170").Value = 3.4
$ws.Range("L170").Value = 3.2
$ws.Range("M170").Value = 2
$ws.Range("N170").Value = 3.4
$ws.Range("O170").Value = 3.1
$ws.Range("P170").Value = -0.25
$ws.Range("Q170").Value = 1.825
$ws.Range("R170").Value = 1.975
$ws.Range("S170").Value = 2.5
$ws.Range("T170").Value = 1.9
$ws.Range("U170").Value = 1.9
$ws.Range("V170").Value = -1
$ws.Range("W170").Value = -1
$ws.Range("X170").Value = 2.1
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = 0.9750000000000001
$ws.Range("AA170").Value = -1
$ws.Range("AB170").Value = 0.8999999999999999

# Row 171
$ws.Range("B171").Value = 6804261
$ws.Range("E171").Value = "Vysocina Jihlava"
$ws.Range("F171").Value = "FC Sellier  Bellot Vlasim"
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = 1
$ws.Range("I171").Value = "D"
$ws.Range("J171").Value = 2.5
$ws.Range("K171").Value = 3.1
$ws.Range("L171").Value = 2.55
$ws.Range("M171").Value = 2.45
$ws.Range("N171").Value = 3.1
$ws.Range("O171").Value = 2.6
$ws.Range("P171").Value = 0
$ws.Range("Q171").Value = 1.8
$ws.Range("R171").Value = 2
$ws.Range("S171").Value = 2.5
$ws.Range("T171").Value = 1.8
$ws.Range("U171").Value = 2
$ws.Range("V171").Value = -1
$ws.Range("W171").Value = 2.1
$ws.Range("X171").Value = -1
$ws.Range("Y171").Value = 0
$ws.Range("Z171").Value = 0
$ws.Range("AA171").Value = -1
$ws.Range("AB171").Value = 1

# Row 172
$ws.Range("B172").Value = 6804263
$ws.Range("E172").Value = "Opava"
$ws.Range("F172").Value = "Sparta Prague B"
$ws.Range("G172").Value = 3
$ws.Range("H172").Value = 2
$ws.Range("I172").Value = "H"
$ws.Range("J172").Value = 2.2
$ws.Range("K172").Value = 3.3
$ws.Range("L172").Value = 2.8
$ws.Range("M172").Value = 1.75
$ws.Range("N172").Value = 3.3
$ws.Range("O172").Value = 4.2
$ws.Range("P172").Value = -0.5
$ws.Range("Q172").Value = 1.825
$ws.Range("R172").Value = 1.975
$ws.Range("S172").Value = 2.25
$ws.Range("T172").Value = 1.925
$ws.Range("U172").Value = 1.875
$ws.Range("V172").Value = 0.75
$ws.Range("W172").Value = -1
$ws.Range("X172").Value = -1
$ws.Range("Y172").Value = 0.825
$ws.Range("Z172").Value = -1
$ws.Range("AA172").Value = 0.925
$ws.Range("AB172").Value = -1

# Row 184
$ws.Range("B184").Value = 6803630
$ws.Range("E184").Value = "FC Brno"
$ws.Range("F184").Value = "Vysocina Jihlava"
$ws.Range("G184").Value = 1
$ws.Range("H184").Value = 0
$ws.Range("I184").Value = "H"
$ws.Range("J184").Value = 1.909
$ws.Range("K184").Value = 3.6
$ws.Range("L184").Value = 3.2
$ws.Range("M184").Value = 1.571
$ws.Range("N184").Value = 3.8
$ws.Range("O184").Value = 4.333
$ws.Range("P184").Value = -0.75
$ws.Range("Q184").Value = 1.775
$ws.Range("R184").Value = 2.025
$ws.Range("S184").Value = 2.75
$ws.Range("T184").Value = 1.775
$ws.Range("U184").Value = 2.025
$ws.Range("V184").Value = 0.571
$ws.Range("W184").Value = -1
$ws.Range("X184").Value = -1
$ws.Range("Y184").Value = 0.3875
$ws.Range("Z184").Value = -0.5
$ws.Range("AA184").Value = -1
$ws.Range("AB184").Value = 1.025

# Row 185
$ws.Range("B185").Value = 6804267
$ws.Range("E185").Value = "FC Sellier  Bellot Vlasim"
$ws.Range("F185").Value = "SK Lisen"
$ws.Range("G185").Value = 4
$ws.Range("H185").Value = 5
$ws.Range("I185").Value = "A"
$ws.Range("J185").Value = 2.15
$ws.Range("K185").Value = 3.1
$ws.Range("L185").Value = 3.1
$ws.Range("M185").Value = 2.1
$ws.Range("N185").Value = 3.1
$ws.Range("O185").Value = 3.3
$ws.Range("P185").Value = -0.25
$ws.Range("Q185").Value = 1.8
$ws.Range("R185").Value = 2
$ws.Range("S185").Value = 2.25
$ws.Range("T185").Value = 1.85
$ws.Range("U185").Value = 1.95
$ws.Range("V185").Value = -1
$ws.Range("W185").Value = -1
$ws.Range("X185").Value = 2.3
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = 1
$ws.Range("AA185").Value = 0.8500000000000001
$ws.Range("AB185").Value = -1

# Row 194
$ws.Range("B194").Value = 6804275
$ws.Range("E194").Value = "MFK Chrudim"
$ws.Range("F194").Value = "FK Pribram"
$ws.Range("G194").Value = 6
$ws.Range("H194").Value = 3
$ws.Range("I194").Value = "H"
$ws.Range("J194").Value = 2.1
$ws.Range("K194").Value = 3.3
$ws.Range("L194").Value = 3
$ws.Range("M194").Value = 2.05
$ws.Range("N194").Value = 3.4
$ws.Range("O194").Value = 3.25
$ws.Range("P194").Value = -0.25
$ws.Range("Q194").Value = 1.8
$ws.Range("R194").Value = 2
$ws.Range("S194").Value = 2.5
$ws.Range("T194").Value = 1.85
$ws.Range("U194").Value = 1.95
$ws.Range("V194").Value = 1.05
$ws.Range("W194").Value = -1
$ws.Range("X194").Value = -1
$ws.Range("Y194").Value = 0.8
$ws.Range("Z194").Value = -1
$ws.Range("AA194").Value = 0.8500000000000001
$ws.Range("AB194").Value = -1

# Row 195
$ws.Range("B195").Value = 6804278
$ws.Range("E195").Value = "MFK Vyskov"
$ws.Range("F195").Value = "Varnsdorf"
$ws.Range("G195").Value = 2
$ws.Range("H195").Value = 1
$ws.Range("I195").Value = "H"
$ws.Range("J195").Value = 1.666
$ws.Range("K195").Value = 3.5
$ws.Range("L195").Value = 4.333
$ws.Range("M195").Value = 1.571
$ws.Range("N195").Value = 3.6
$ws.Range("O195").Value = 5
$ws.Range("P195").Value = -0.75
$ws.Range("Q195").Value = 1.8
$ws.Range("R195").Value = 2
$ws.Range("S195").Value = 2.75
$ws.Range("T195").Value = 1.825
$ws.Range("U195").Value = 1.975
$ws.Range("V195").Value = 0.571
$ws.Range("W195").Value = -1
$ws.Range("X195").Value = -1
$ws.Range("Y195").Value = 0.4
$ws.Range("Z195").Value = -0.5
$ws.Range("AA195").Value = 0.4125
$ws.Range("AB195").Value = -0.5

# Row 210
$ws.Range("B210").Value = 6804287
$ws.Range("E210").Value = "Sparta Prague B"
$ws.Range("F210").Value = "SK Lisen"
$ws.Range("G210").Value = 1
$ws.Range("H210").Value = 3
$ws.Range("I210").Value = "A"
$ws.Range("J210").Value = 1.8
$ws.Range("K210").Value = 3.5
$ws.Range("L210").Value = 3.6
$ws.Range("M210").Value = 1.85
$ws.Range("N210").Value = 3.6
$ws.Range("O210").Value = 3.8
$ws.Range("P210").Value = -0.5
$ws.Range("Q210").Value = 1.85
$ws.Range("R210").Value = 1.95
$ws.Range("S210").Value = 2.75
$ws.Range("T210").Value = 1.975
$ws.Range("U210").Value = 1.825
$ws.Range("V210").Value = -1
$ws.Range("W210").Value = -1
$ws.Range("X210").Value = 2.8
$ws.Range("Y210").Value = -1
$ws.Range("Z210").Value = 0.95
$ws.Range("AA210").Value = 0.9750000000000001
$ws.Range("AB210").Value = -1

# Row 211
$ws.Range("B211").Value = 6804285
$ws.Range("E211").Value = "SK Prostejov"
$ws.Range("F211").Value = "FK Pribram"
$ws.Range("G211").Value = 2
$ws.Range("H211").Value = 1
$ws.Range("I211").Value = "H"
$ws.Range("J211").Value = 2
$ws.Range("K211").Value = 3.75
$ws.Range("L211").Value = 2.9
$ws.Range("M211").Value = 1.8
$ws.Range("N211").Value = 3.8
$ws.Range("O211").Value = 3.4
$ws.Range("P211").Value = -0.5
$ws.Range("Q211").Value = 1.85
$ws.Range("R211").Value = 1.95
$ws.Range("S211").Value = 2.75
$ws.Range("T211").Value = 1.9
$ws.Range("U211").Value = 1.9
$ws.Range("V211").Value = 0.8
$ws.Range("W211").Value = -1
$ws.Range("X211").Value = -1
$ws.Range("Y211").Value = 0.8500000000000001
$ws.Range("Z211").Value = -1
$ws.Range("AA211").Value = 0.45
$ws.Range("AB211").Value = -0.5

# Row 213
$ws.Range("B213").Value = 6803349
$ws.Range("E213").Value = "Varnsdorf"
$ws.Range("F213").Value = "Hanacka Slavia Kromeriz"
$ws.Range("G213").Value = 3
$ws.Range("H213").Value = 2
$ws.Range("I213").Value = "H"
$ws.Range("J213").Value = 1.571
$ws.Range("K213").Value = 3.8
$ws.Range("L213").Value = 4.5
$ws.Range("M213").Value = 1.533
$ws.Range("N213").Value = 4
$ws.Range("O213").Value = 4.75
$ws.Range("P213").Value = -1
$ws.Range("Q213").Value = 1.9
$ws.Range("R213").Value = 1.9
$ws.Range("S213").Value = 3.25
$ws.Range("T213").Value = 2
$ws.Range("U213").Value = 1.8
$ws.Range("V213").Value = 0.5329999999999999
$ws.Range("W213").Value = -1
$ws.Range("X213").Value = -1
$ws.Range("Y213").Value = 0
$ws.Range("Z213").Value = 0
$ws.Range("AA213").Value = 1
$ws.Range("AB213").Value = -1

# Row 214
$ws.Range("B214").Value = 6804286
$ws.Range("E214").Value = "FC Silon Taborsko"
$ws.Range("F214").Value = "Opava"
$ws.Range("G214").Value = 2
$ws.Range("H214").Value = 1
$ws.Range("I214").Value = "H"
$ws.Range("J214").Value = 1.727
$ws.Range("K214").Value = 3.5
$ws.Range("L214").Value = 4
$ws.Range("M214").Value = 1.8
$ws.Range("N214").Value = 3.3
$ws.Range("O214").Value = 3.8
$ws.Range("P214").Value = -0.5
$ws.Range("Q214").Value = 1.875
$ws.Range("R214").Value = 1.925
$ws.Range("S214").Value = 2.25
$ws.Range("T214").Value = 1.875
$ws.Range("U214").Value = 1.925
$ws.Range("V214").Value = 0.8
$ws.Range("W214").Value = -1
$ws.Range("X214").Value = -1
$ws.Range("Y214").Value = 0.875
$ws.Range("Z214").Value = -1
$ws.Range("AA214").Value = 0.875
$ws.Range("AB214").Value = -1

# Row 218
$ws.Range("B218").Value = 6803350
$ws.Range("E218").Value = "Hanacka Slavia Kromeriz"
$ws.Range("F218").Value = "MFK Chrudim"
$ws.Range("G218").Value = 1
$ws.Range("H218").Value = 1
$ws.Range("I218").Value = "D"
$ws.Range("J218").Value = 2.5
$ws.Range("K218").Value = 3.25
$ws.Range("L218").Value = 2.45
$ws.Range("M218").Value = 2.8
$ws.Range("N218").Value = 3.4
$ws.Range("O218").Value = 2.15
$ws.Range("P218").Value = 0.25
$ws.Range("Q218").Value = 1.825
$ws.Range("R218").Value = 1.975
$ws.Range("S218").Value = 2.75
$ws.Range("T218").Value = 1.95
$ws.Range("U218").Value = 1.85
$ws.Range("V218").Value = -1
$ws.Range("W218").Value = 2.4
$ws.Range("X218").Value = -1
$ws.Range("Y218").Value = 0.4125
$ws.Range("Z218").Value = -0.5
$ws.Range("AA218").Value = -1
$ws.Range("AB218").Value = 0.8500000000000001

# Row 219
$ws.Range("B219").Value = 6804296
$ws.Range("E219").Value = "FK Pribram"
$ws.Range("F219").Value = "Vysocina Jihlava"
$ws.Range("G219").Value = 1
$ws.Range("H219").Value = 1
$ws.Range("I219").Value = "D"
$ws.Range("J219").Value = 2.05
$ws.Range("K219").Value = 3.5
$ws.Range("L219").Value = 2.9
$ws.Range("M219").Value = 1.85
$ws.Range("N219").Value = 3.5
$ws.Range("O219").Value = 3.4
$ws.Range("P219").Value = -0.5
$ws.Range("Q219").Value = 1.925
$ws.Range("R219").Value = 1.875
$ws.Range("S219").Value = 2.5
$ws.Range("T219").Value = 1.9
$ws.Range("U219").Value = 1.9
$ws.Range("V219").Value = -1
$ws.Range("W219").Value = 2.5
$ws.Range("X219").Value = -1
$ws.Range("Y219").Value = -1
$ws.Range("Z219").Value = 0.875
$ws.Range("AA219").Value = -1
$ws.Range("AB219").Value = 0.8999999999999999

# Row 223
$ws.Range("B223").Value = 6804294
$ws.Range("E223").Value = "Opava"
$ws.Range("F223").Value = "MFK Vyskov"
$ws.Range("G223").Value = 3
$ws.Range("H223").Value = 1
$ws.Range("I223").Value = "H"
$ws.Range("J223").Value = 2.15
$ws.Range("K223").Value = 3.1
$ws.Range("L223").Value = 3.1
$ws.Range("M223").Value = 1.85
$ws.Range("N223").Value = 3.25
$ws.Range("O223").Value = 3.8
$ws.Range("P223").Value = -0.5
$ws.Range("Q223").Value = 1.925
$ws.Range("R223").Value = 1.875
$ws.Range("S223").Value = 2.5
$ws.Range("T223").Value = 2
$ws.Range("U223").Value = 1.8
$ws.Range("V223").Value = 0.8500000000000001
$ws.Range("W223").Value = -1
$ws.Range("X223").Value = -1
$ws.Range("Y223").Value = 0.925
$ws.Range("Z223").Value = -1
$ws.Range("AA223").Value = 1
$ws.Range("AB223").Value = -1

# Row 224
$ws.Range("B224").Value = 6804293
$ws.Range("E224").Value = "Dukla Praha"
$ws.Range("F224").Value = "SK Prostejov"
$ws.Range("G224").Value = 3
$ws.Range("H224").Value = 2
$ws.Range("I224").Value = "H"
$ws.Range("J224").Value = 1.4
$ws.Range("K224").Value = 4.5
$ws.Range("L224").Value = 5.25
$ws.Range("M224").Value = 1.533
$ws.Range("N224").Value = 4
$ws.Range("O224").Value = 4.75
$ws.Range("P224").Value = -1
$ws.Range("Q224").Value = 1.9
$ws.Range("R224").Value = 1.9
$ws.Range("S224").Value = 3
$ws.Range("T224").Value = 1.95
$ws.Range("U224").Value = 1.85
$ws.Range("V224").Value = 0.5329999999999999
$ws.Range("W224").Value = -1
$ws.Range("X224").Value = -1
$ws.Range("Y224").Value = 0
$ws.Range("Z224").Value = 0
$ws.Range("AA224").Value = 0.95
$ws.Range("AB224").Value = -1
